$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get reshuffled across rows 2-18: D, M, N, O, P, S
$cols = @("D", "M", "N", "O", "P", "S")

# Snapshot the current (before) values for each relevant column, rows 2-18
$orig = @{}
foreach ($col in $cols) {
    $orig[$col] = @{}
    for ($r = 2; $r -le 18; $r++) {
        $orig[$col][$r] = $ws.Range("$col$r").Value2
    }
}

# Mapping: new row -> source (old) row, i.e. row $r takes the old values of row $perm[$r]
$perm = @{
    2  = 11
    3  = 2
    4  = 5
    5  = 3
    6  = 4
    7  = 13
    8  = 18
    9  = 14
    10 = 12
    11 = 8
    12 = 9
    13 = 17
    14 = 10
    15 = 7
    16 = 6
    17 = 16
    18 = 15
}

foreach ($r in $perm.Keys) {
    $src = $perm[$r]
    foreach ($col in $cols) {
        $ws.Range("$col$r").Value = $orig[$col][$src]
    }
}
